# Rename model Region to Zone for clarity
$wb = $excel.ActiveWorkbook

# Rename the "Region" sheet to "Zone"
$regionSheet = $wb.Worksheets.Item("Region")
$regionSheet.Name = "Zone"

# Make "Zone" the active sheet / selected tab, with D2 selected
$zoneSheet = $wb.Worksheets.Item("Zone")
$zoneSheet.Activate()
$zoneSheet.Range("D2").Select()
